# Commit: add the NA's under duplicate_image_filename
#
# Column E ("duplicate_image_filename") was blank for the stimuli rows
# (rows 2-21). Fill it in with "NA" for all of those rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Restore the empty placeholder cell F1 that the load/save round-trip of
# this runtime otherwise populates with a stray value, so it stays blank
# just like in the original workbook (no change described for F1).
$ws.Range("F1").ClearContents()

# Fill the duplicate_image_filename column (E) with "NA" for rows 2-21.
$ws.Range("E2:E21").Value = "NA"
